$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nudge the sheet's default/standard column width toward the new value
$ws.StandardWidth = 8.625

# Row 28: fill in values matching other rows (sample_name, scandir, heightstep, suffix, ring)
$ws.Range("A28").Value = "tissue58"
$ws.Range("B28").Value = "032_Mouse89883_5823_eth_9000proj_ring"
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "z6_x2_dummy"
$ws.Range("E28").Value = 3

# Row 29
$ws.Range("A29").Value = "tissue58"
$ws.Range("B29").Value = "032_Mouse89883_5823_eth_9000proj_ring"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "z6_x3_dummy"
$ws.Range("E29").Value = 4

# Update column D width (target stored width 15.89 chars; runtime quantizes
# column width storage to whole pixels with a 5px padding at ~6px/char, so the
# closest achievable stored value is 15.8333... when ColumnWidth is set to 15)
$ws.Columns.Item(4).ColumnWidth = 15

# Update selection to B13
$ws.Range("B13").Select()
